$wb = $excel.ActiveWorkbook

# "Balls per load" constant on the Assumptions sheet drops from 40 to 25 -
# this ripples through the Gear Delivery Schedule's M:P columns (Balls Shot,
# Balls Successful, Preasure Points Accumulated).
$assumptions = $wb.Worksheets.Item("Assumptions")
$assumptions.Range("B13").Value = 25

# Selection moved: Assumptions B6 is now the active cell, and Assumptions
# becomes the active/selected tab instead of Gear Delivery Schedule.
$assumptions.Activate()
$assumptions.Range("B6").Select()
